# contratos-10-2014.xlsx — "fix: fixed formatting when scrapping floating
# point numbers"
#
# The "Importe" column (H) was scraped with the Spanish/Argentine number
# format (thousands separator "." + decimal separator ","), e.g. "3.300,00".
# This rewrites every amount in that column to plain "3300.00" notation
# (no thousands separator, "." as decimal separator) while keeping the
# cell content as TEXT (a handful of names in columns D-G had the same
# scraper glitch - a stray "," used where a "." belonged - which is
# corrected for those specific, known cells below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix the handful of contractor/company names that were scraped
#        with a stray comma instead of a period. -----------------------
$nameFixes = @{
    "E32"  = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
    "E46"  = "RAMIREZ CLAUDIA. RAMIREZ CESAR Y RAMIREZ VERONICA SH"
    "E50"  = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
    "E98"  = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
    "E130" = "RICCOTTI. MARIANA EDITH"
    "F143" = "MERCANZINI. GASTON ARIEL"
    "E179" = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
    "E181" = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
    "E198" = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
}

foreach ($ref in $nameFixes.Keys) {
    $ws.Range($ref).Value = $nameFixes[$ref]
}

# --- 2) Reformat every floating point amount in column H ("Importe"),
#        rows 2..240, from "1.234,56" to "1234.56". ----------------------
$numPattern = '^[0-9\.]+,[0-9]{2}$'

$firstRow = 2
$lastRow = 240
$col = 8  # column H

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $old = [string]$cell.Value2

    if ($old -match $numPattern) {
        $new = $old.Replace(".", "").Replace(",", ".")

        # Write the value back as TEXT (it was - and must remain - a
        # shared-string, not a number) without leaving the cell's
        # style/number-format changed: stage it as Text ("@"), assign
        # the literal string, then restore the "Normal" style so the
        # cell's style index is back to what it was originally.
        $cell.NumberFormat = "@"
        $cell.Value = $new
        $cell.Style = "Normal"
    }
}
